$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# Change 1: drop the stray "_GoBack" bookmark pair that trailed the
# "coal (commodity = coal)" bullet. Locate that paragraph by its text and
# rewrite the paragraph's Range with identical formatting/text but no
# bookmark markers.
# ---------------------------------------------------------------------------
$coalXml = @"
<w:p xmlns:w="$wNs"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>coal (commodity = coal)</w:t></w:r></w:p>
"@

$coalFound = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq "coal (commodity = coal)") {
        [void]$para.Range.InsertXML($coalXml)
        $coalFound = $true
        break
    }
}
if (-not $coalFound) {
    Write-Output "WARNING: 'coal (commodity = coal)' paragraph not found"
}

# ---------------------------------------------------------------------------
# Change 2: the document used to end with one empty "ListParagraph" bullet.
# Turn it into a bold "output" heading, then append the new Q&A content that
# follows it in the source doc: year_vtg/year_act question, capacity_factor
# heading + question, ref_activity heading + question, relation_activity
# heading + question (the last of which re-introduces the "_GoBack"
# bookmark, now in its new spot).
# ---------------------------------------------------------------------------
$outputPara = @"
<w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>output</w:t></w:r></w:p>
"@

$yearPara = @"
<w:p xmlns:w="$wNs"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">How are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>year_vtg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>year_act</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> set?</w:t></w:r></w:p>
"@

$capacityFactorHeading = @"
<w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>capacity_factor</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@

$capacityFactorQuestion = @"
<w:p xmlns:w="$wNs"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Is it ok that CF = 1 for all trade technologies?</w:t></w:r></w:p>
"@

$refActivityHeading = @"
<w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>ref_activity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@

$refActivityQuestion = @"
<w:p xmlns:w="$wNs"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">How should we build </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ref_activity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>?</w:t></w:r></w:p>
"@

$relationActivityHeading = @"
<w:p xmlns:w="$wNs"><w:pPr><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>relation_activity</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@

$relationActivityQuestion = @"
<w:p xmlns:w="$wNs"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">What is </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>PE_total_engineering</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>?</w:t></w:r></w:p>
"@

$tailXml = $outputPara + $yearPara + $capacityFactorHeading + $capacityFactorQuestion + `
           $refActivityHeading + $refActivityQuestion + $relationActivityHeading + $relationActivityQuestion

$insertionPoint = $d.Paragraphs.Last.Range
$insertionPoint.Collapse(0)
[void]$insertionPoint.InsertXML($tailXml)

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
